$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 98.166664
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H40").Value = 2069.4443
$ws.Range("I40").Value = 1589.2858
$ws.Range("K40").Value = 1589.2858
$ws.Range("M40").Value = -1414.2858
$ws.Range("H70").Value = 2025.091
$ws.Range("I70").Value = 1664.8334
$ws.Range("J70").Value = 2457.4
$ws.Range("K70").Value = 4994.5002
$ws.Range("L70").Value = 7372.200000000001
$ws.Range("M70").Value = -4724.5002
$ws.Range("N70").Value = -7912.200000000001
$ws.Range("H73").Value = 2025.091
$ws.Range("I73").Value = 1664.8334
$ws.Range("J73").Value = 2457.4
$ws.Range("K73").Value = 4994.5002
$ws.Range("L73").Value = 7372.200000000001
$ws.Range("M73").Value = -4058.5002
$ws.Range("N73").Value = -9244.200000000001
$ws.Range("H76").Value = 5858.76
$ws.Range("I76").Value = 5384.15
$ws.Range("K76").Value = 5384.15
$ws.Range("M76").Value = -5069.15
$ws.Range("H79").Value = 5858.76
$ws.Range("I79").Value = 5384.15
$ws.Range("K79").Value = 5384.15
$ws.Range("M79").Value = -4292.15
$ws.Range("H80").Value = 886
$ws.Range("I80").Value = 543.75
$ws.Range("J80").Value = 1159.8
$ws.Range("K80").Value = 1631.25
$ws.Range("L80").Value = 3479.4
$ws.Range("M80").Value = -633.25
$ws.Range("N80").Value = -5475.4
$ws.Range("H83").Value = 886
$ws.Range("I83").Value = 543.75
$ws.Range("J83").Value = 1159.8
$ws.Range("K83").Value = 4893.75
$ws.Range("L83").Value = 10438.2
$ws.Range("M83").Value = 98.25
$ws.Range("N83").Value = -20422.2
$ws.Range("H86").Value = 2860522
$ws.Range("J86").Value = 7145206.5
$ws.Range("L86").Value = 7145206.5
$ws.Range("N86").Value = -7147452.5
$ws.Range("H88").Value = 479070.25
$ws.Range("I88").Value = 2334.8572
$ws.Range("K88").Value = 2334.8572
$ws.Range("M88").Value = -1928.8572
$ws.Range("H89").Value = 2860522
$ws.Range("J89").Value = 7145206.5
$ws.Range("L89").Value = 35726032.5
$ws.Range("N89").Value = -35737264.5
$ws.Range("H91").Value = 479070.25
$ws.Range("I91").Value = 2334.8572
$ws.Range("K91").Value = 2334.8572
$ws.Range("M91").Value = -930.8571999999999
$ws.Range("H112").Value = 2515.1667
$ws.Range("J112").Value = 2838.2
$ws.Range("L112").Value = 8514.599999999999
$ws.Range("N112").Value = -10730.6
$ws.Range("H138").Value = 2428.5576
$ws.Range("J138").Value = 2857.5898
$ws.Range("L138").Value = 8572.769400000001
$ws.Range("N138").Value = -18852.7694
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3391.5625
$ws.Range("J88").Value = 3530.25
$ws.Range("L88").Value = 3530.25
$ws.Range("N88").Value = -4342.25
$ws.Range("H91").Value = 3391.5625
$ws.Range("J91").Value = 3530.25
$ws.Range("L91").Value = 3530.25
$ws.Range("N91").Value = -6338.25
$ws.Range("H132").Value = 8350.654
$ws.Range("I132").Value = 4432.4683
$ws.Range("K132").Value = 13297.4049
$ws.Range("M132").Value = -10767.4049
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4899.8335
$ws.Range("I86").Value = 2999.6667
$ws.Range("K86").Value = 2999.6667
$ws.Range("M86").Value = -1876.6667
$ws.Range("H89").Value = 4899.8335
$ws.Range("I89").Value = 2999.6667
$ws.Range("K89").Value = 14998.3335
$ws.Range("M89").Value = -9382.333500000001
$ws.Range("H107").Value = 8226.5
$ws.Range("I107").Value = 8226.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 8226.5
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2281.0588
$ws.Range("I31").Value = 2434.3635
$ws.Range("K31").Value = 2434.3635
$ws.Range("M31").Value = -2139.3635
$ws.Range("H34").Value = 2281.0588
$ws.Range("I34").Value = 2434.3635
$ws.Range("K34").Value = 2434.3635
$ws.Range("M34").Value = -2232.3635
$ws.Range("H99").Value = 4135.6875
$ws.Range("I99").Value = 3696.7778
$ws.Range("J99").Value = 4700
$ws.Range("K99").Value = 3696.7778
$ws.Range("L99").Value = 4700
$ws.Range("M99").Value = -2198.7778
$ws.Range("N99").Value = -7696
$ws.Range("H107").Value = 891.5263
$ws.Range("I107").Value = 476.125
$ws.Range("K107").Value = 476.125
$ws.Range("M107").Value = 1443.875
$ws.Range("H126").Value = 4135.6875
$ws.Range("I126").Value = 3696.7778
$ws.Range("J126").Value = 4700
$ws.Range("K126").Value = 11090.3334
$ws.Range("L126").Value = 14100
$ws.Range("M126").Value = -8620.3334
$ws.Range("N126").Value = -19040
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 27777934
$ws.Range("I2").Value = 143.71428
$ws.Range("K2").Value = 862.28568
$ws.Range("M2").Value = -749.28568
$ws.Range("H37").Value = 84982
$ws.Range("J37").Value = 84982
$ws.Range("L37").Value = 254946
$ws.Range("N37").Value = -255170
$ws.Range("H38").Value = 152.8077
$ws.Range("J38").Value = 166.55556
$ws.Range("L38").Value = 499.66668
$ws.Range("N38").Value = -1193.66668
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3912.9524
$ws.Range("I80").Value = 3912.9524
$ws.Range("K80").Value = 3912.9524
$ws.Range("M80").Value = -2914.9524
$ws.Range("H83").Value = 3912.9524
$ws.Range("I83").Value = 3912.9524
$ws.Range("K83").Value = 19564.762
$ws.Range("M83").Value = -14572.762
$ws.Range("H132").Value = 1986.5
$ws.Range("I132").Value = 1986.5
$ws.Range("K132").Value = 5959.5
$ws.Range("M132").Value = -3429.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 50000
$ws.Range("J63").Value = 50000
$ws.Range("L63").Value = 50000
$ws.Range("N63").Value = -51498
$ws.Range("H66").Value = 50000
$ws.Range("J66").Value = 50000
$ws.Range("L66").Value = 150000
$ws.Range("N66").Value = -157488
$ws.Range("H82").Value = 2932.5
$ws.Range("I82").Value = 806.7778
$ws.Range("J82").Value = 5665.5713
$ws.Range("K82").Value = 806.7778
$ws.Range("L82").Value = 5665.5713
$ws.Range("M82").Value = -445.7778
$ws.Range("N82").Value = -6387.5713
$ws.Range("H85").Value = 2932.5
$ws.Range("I85").Value = 806.7778
$ws.Range("J85").Value = 5665.5713
$ws.Range("K85").Value = 806.7778
$ws.Range("L85").Value = 5665.5713
$ws.Range("M85").Value = 441.2222
$ws.Range("N85").Value = -8161.5713
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I68").Value = 2000
$ws.Range("K68").Value = 2000
$ws.Range("M68").Value = -1189
$ws.Range("I71").Value = 2000
$ws.Range("K71").Value = 6000
$ws.Range("M71").Value = -1944
$ws.Range("H81").Value = 3115.8096
$ws.Range("I81").Value = 1978.6666
$ws.Range("J81").Value = 3968.6667
$ws.Range("K81").Value = 3957.3332
$ws.Range("L81").Value = 7937.3334
$ws.Range("M81").Value = -2896.3332
$ws.Range("N81").Value = -10059.3334
$ws.Range("H84").Value = 3115.8096
$ws.Range("I84").Value = 1978.6666
$ws.Range("J84").Value = 3968.6667
$ws.Range("K84").Value = 19786.666
$ws.Range("L84").Value = 39686.667
$ws.Range("M84").Value = -14482.666
$ws.Range("N84").Value = -50294.667
$ws.Range("H126").Value = 2167
$ws.Range("J126").Value = 2250.5
$ws.Range("L126").Value = 6751.5
$ws.Range("N126").Value = -11691.5
$ws.Range("H132").Value = 1450.3429
$ws.Range("I132").Value = 1450.3429
$ws.Range("K132").Value = 4351.028700000001
$ws.Range("M132").Value = -1821.028700000001
Write-Host "Applied all Zodiark_Profits updates."
